# "Generate Report for Archive"
#
# The status text "Ready for handoff" is replaced with "In Translation"
# everywhere it appears (it is a single shared string used by the
# Overview sheet's E2/F2 status cells and by the zh-cn / de-de sheets'
# C2 status cells). Because the shorter replacement text no longer needs
# as much horizontal space, the status column on each sheet is narrowed
# to match (mirrors an autofit-after-edit pass).

$wb = $excel.ActiveWorkbook

$oldStatus = "Ready for handoff"
$newStatus = "In Translation"

# --- Overview sheet: status columns are E (zh-cn) and F (de-de) ---
$wsOverview = $wb.Worksheets.Item("Overview")
if ($wsOverview.Range("E2").Value2 -eq $oldStatus) {
    $wsOverview.Range("E2").Value2 = $newStatus
}
if ($wsOverview.Range("F2").Value2 -eq $oldStatus) {
    $wsOverview.Range("F2").Value2 = $newStatus
}
$wsOverview.Columns.Item(5).ColumnWidth = 12.5
$wsOverview.Columns.Item(6).ColumnWidth = 12.5

# --- zh-cn sheet: status column is C ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
if ($wsZhCn.Range("C2").Value2 -eq $oldStatus) {
    $wsZhCn.Range("C2").Value2 = $newStatus
}
$wsZhCn.Columns.Item(3).ColumnWidth = 12.5

# --- de-de sheet: status column is C ---
$wsDeDe = $wb.Worksheets.Item("de-de")
if ($wsDeDe.Range("C2").Value2 -eq $oldStatus) {
    $wsDeDe.Range("C2").Value2 = $newStatus
}
$wsDeDe.Columns.Item(3).ColumnWidth = 12.5
